# Slide 2 ("Where are we?"): update the "http://go/..." link text/URL.
#
# Original paragraph 1 of the content placeholder reads:
#   "Have a look at http://go/5hciuvsg "
# and is made of two runs: "Have a look at " and "http://go/5hciuvsg ".
#
# Target state (per the authored diff) splits that into six runs:
#   "Have a look " | "at " | "http" | "://" | "go/x2std3ms" | " "
# with the last three ("http", "://", "go/x2std3ms") carrying a
# hyperlink to the (updated) go-link target "http://go/x2std3ms".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$para1 = $tr.Paragraphs(1, 1)

$newUrl = "http://go/x2std3ms"

# --- Split "Have a look at " -> "Have a look " + "at " -----------------
# Characters(13,3) is the "at " substring (1-based); reassigning its own
# text forces a clean run boundary there without touching the preceding
# "Have a look " run.
$atPiece = $para1.Characters(13, 3)
$atPiece.Text = "at "

# --- Turn "http://go/5hciuvsg " into a hyperlinked, re-targeted URL ----
# First carve out the three pieces that will carry the hyperlink
# ("http", "://", "go/5hciuvsg") plus the trailing space that stays
# plain text, then point them at the new go-link.
$httpPiece  = $para1.Characters(16, 4)   # "http"
$slashPiece = $para1.Characters(20, 3)   # "://"
$goPiece    = $para1.Characters(23, 11)  # "go/5hciuvsg"

$httpPiece.ActionSettings.Item(1).Hyperlink.Address  = $newUrl
$slashPiece.ActionSettings.Item(1).Hyperlink.Address = $newUrl
$goPiece.ActionSettings.Item(1).Hyperlink.Address    = $newUrl

# Update the visible slug from "5hciuvsg" to "x2std3ms".
$goPiece.Text = "go/x2std3ms"

Write-Output $para1.Text
